$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - all subsequent rows shift up by one
$ws.Rows.Item(26).Delete()

# After the above deletion, the row that used to be "SC 92" (was row 28) is now
# at row 27. Delete it too, shifting everything up by one more row.
$ws.Rows.Item(27).Delete()

# At this point the sheet data (rows 26-33) reads:
#   26 SC 5    -20.2 10.8  -13.8 -5    17.38
#   27 SC 101  -20.4 (missing) -14.6 -10   17
#   28 SC 105  -19.6 11.1  -13.7 -5.9  17.44
#   29 SC 119  -19.5 11.2  -13   -6.8  18.06
#   30 SC 120  -19.7 (missing) -13.6 -5.7  16.89
#   31 SC 132  -18.8 15.3  -13.7 -8.1  17.18
#   32 SC 193  -19.9 10.5  -14.7 -6.4  17.39
#   33 SC 232  -19.5 10.4  -14.1 -10.7 17.53
#
# The target missing-data pattern requires updating column C ("B" header) values
# for several rows: fill in previously-missing cells and blank out others.

$ws.Range("C27").Value = 10
$ws.Range("C28").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("C30").Value = 11.4
$ws.Range("C32").ClearContents()
